# Auto-generated Excel COM-interop script to update Leve profit data
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 267.6  # H2
$ws.Cells.Item(2, 9).Value = 272  # I2
$ws.Cells.Item(2, 11).Value = 272  # K2
$ws.Cells.Item(2, 13).Value = -159  # M2
$ws.Cells.Item(33, 8).Value = 936.0645  # H33
$ws.Cells.Item(33, 9).Value = 779.087  # I33
$ws.Cells.Item(33, 11).Value = 779.087  # K33
$ws.Cells.Item(33, 13).Value = -550.087  # M33
$ws.Cells.Item(138, 8).Value = 2111.4922  # H138
$ws.Cells.Item(138, 9).Value = 2368.5  # I138
$ws.Cells.Item(138, 11).Value = 7105.5  # K138
$ws.Cells.Item(138, 13).Value = -1965.5  # M138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 472717.06  # H32
$ws.Cells.Item(32, 9).Value = 617522.4  # I32
$ws.Cells.Item(32, 11).Value = 617522.4  # K32
$ws.Cells.Item(32, 13).Value = -617235.4  # M32
$ws.Cells.Item(61, 8).Value = 4112.5  # H61
$ws.Cells.Item(61, 9).Value = 2966.6667  # I61
$ws.Cells.Item(61, 10).Value = 4800  # J61
$ws.Cells.Item(61, 11).Value = 2966.6667  # K61
$ws.Cells.Item(61, 12).Value = 4800  # L61
$ws.Cells.Item(61, 13).Value = -2754.6667  # M61
$ws.Cells.Item(61, 14).Value = -5224  # N61
$ws.Cells.Item(74, 8).Value = 863.8182  # H74
$ws.Cells.Item(74, 9).Value = 810.3125  # I74
$ws.Cells.Item(74, 10).Value = 1006.5  # J74
$ws.Cells.Item(74, 11).Value = 810.3125  # K74
$ws.Cells.Item(74, 12).Value = 1006.5  # L74
$ws.Cells.Item(74, 13).Value = 63.6875  # M74
$ws.Cells.Item(74, 14).Value = -2754.5  # N74
$ws.Cells.Item(77, 8).Value = 863.8182  # H77
$ws.Cells.Item(77, 9).Value = 810.3125  # I77
$ws.Cells.Item(77, 10).Value = 1006.5  # J77
$ws.Cells.Item(77, 11).Value = 4051.5625  # K77
$ws.Cells.Item(77, 12).Value = 5032.5  # L77
$ws.Cells.Item(77, 13).Value = 316.4375  # M77
$ws.Cells.Item(77, 14).Value = -13768.5  # N77
$ws.Cells.Item(136, 8).Value = 4112.5  # H136
$ws.Cells.Item(136, 9).Value = 2966.6667  # I136
$ws.Cells.Item(136, 10).Value = 4800  # J136
$ws.Cells.Item(136, 11).Value = 8900.000100000001  # K136
$ws.Cells.Item(136, 12).Value = 14400  # L136
$ws.Cells.Item(136, 13).Value = -6350.000100000001  # M136
$ws.Cells.Item(136, 14).Value = -19500  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2273.2856  # H107
$ws.Cells.Item(107, 9).Value = 1000  # I107
$ws.Cells.Item(107, 10).Value = 2485.5  # J107
$ws.Cells.Item(107, 11).Value = 1000  # K107
$ws.Cells.Item(107, 12).Value = 2485.5  # L107
$ws.Cells.Item(107, 13).Value = 920  # M107
$ws.Cells.Item(107, 14).Value = -6325.5  # N107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 26334.334  # H4
$ws.Cells.Item(4, 9).Value = 5000.5  # I4
$ws.Cells.Item(4, 10).Value = 69002  # J4
$ws.Cells.Item(4, 11).Value = 5000.5  # K4
$ws.Cells.Item(4, 12).Value = 69002  # L4
$ws.Cells.Item(4, 13).Value = -4888.5  # M4
$ws.Cells.Item(4, 14).Value = -69226  # N4
$ws.Cells.Item(7, 8).Value = 59.166668  # H7
$ws.Cells.Item(7, 9).Value = 41.666668  # I7
$ws.Cells.Item(7, 10).Value = 76.666664  # J7
$ws.Cells.Item(7, 11).Value = 41.666668  # K7
$ws.Cells.Item(7, 12).Value = 76.666664  # L7
$ws.Cells.Item(7, 13).Value = 71.333332  # M7
$ws.Cells.Item(7, 14).Value = -302.666664  # N7
$ws.Cells.Item(10, 8).Value = 25000  # H10
$ws.Cells.Item(10, 9).Value = 0  # I10
$ws.Cells.Item(10, 11).Value = 0  # K10
$ws.Cells.Item(10, 13).ClearContents()  # M10
$ws.Cells.Item(31, 8).Value = 2434.7646  # H31
$ws.Cells.Item(31, 9).Value = 981.2174  # I31
$ws.Cells.Item(31, 10).Value = 5474  # J31
$ws.Cells.Item(31, 11).Value = 981.2174  # K31
$ws.Cells.Item(31, 12).Value = 5474  # L31
$ws.Cells.Item(31, 13).Value = -686.2174  # M31
$ws.Cells.Item(31, 14).Value = -6064  # N31
$ws.Cells.Item(34, 8).Value = 2434.7646  # H34
$ws.Cells.Item(34, 9).Value = 981.2174  # I34
$ws.Cells.Item(34, 10).Value = 5474  # J34
$ws.Cells.Item(34, 11).Value = 981.2174  # K34
$ws.Cells.Item(34, 12).Value = 5474  # L34
$ws.Cells.Item(34, 13).Value = -779.2174  # M34
$ws.Cells.Item(34, 14).Value = -5878  # N34
$ws.Cells.Item(58, 8).Value = 1698.8572  # H58
$ws.Cells.Item(58, 9).Value = 826.6667  # I58
$ws.Cells.Item(58, 10).Value = 2353  # J58
$ws.Cells.Item(58, 11).Value = 826.6667  # K58
$ws.Cells.Item(58, 12).Value = 2353  # L58
$ws.Cells.Item(58, 13).Value = -623.6667  # M58
$ws.Cells.Item(58, 14).Value = -2759  # N58
$ws.Cells.Item(132, 8).Value = 11112913  # H132
$ws.Cells.Item(132, 9).Value = 1051  # I132
$ws.Cells.Item(132, 10).Value = 33336636  # J132
$ws.Cells.Item(132, 11).Value = 3153  # K132
$ws.Cells.Item(132, 12).Value = 100009908  # L132
$ws.Cells.Item(132, 13).Value = -623  # M132
$ws.Cells.Item(132, 14).Value = -100014968  # N132
$ws.Cells.Item(136, 8).Value = 1698.8572  # H136
$ws.Cells.Item(136, 9).Value = 826.6667  # I136
$ws.Cells.Item(136, 10).Value = 2353  # J136
$ws.Cells.Item(136, 11).Value = 2480.0001  # K136
$ws.Cells.Item(136, 12).Value = 7059  # L136
$ws.Cells.Item(136, 13).Value = 69.9998999999998  # M136
$ws.Cells.Item(136, 14).Value = -12159  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 2859.1428  # H60
$ws.Cells.Item(60, 9).Value = 420  # I60
$ws.Cells.Item(60, 10).Value = 3136.318  # J60
$ws.Cells.Item(60, 11).Value = 1260  # K60
$ws.Cells.Item(60, 12).Value = 9408.954000000002  # L60
$ws.Cells.Item(60, 13).Value = -1009  # M60
$ws.Cells.Item(60, 14).Value = -9910.954000000002  # N60
$ws.Cells.Item(115, 8).Value = 4052.25  # H115
$ws.Cells.Item(115, 10).Value = 4396.636  # J115
$ws.Cells.Item(115, 12).Value = 13189.908  # L115
$ws.Cells.Item(115, 14).Value = -15539.908  # N115

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 27057.666  # H34
$ws.Cells.Item(34, 10).Value = 27057.666  # J34
$ws.Cells.Item(34, 12).Value = 27057.666  # L34
$ws.Cells.Item(34, 14).Value = -27593.666  # N34
$ws.Cells.Item(43, 8).Value = 12795.4  # H43
$ws.Cells.Item(43, 9).Value = 1990  # I43
$ws.Cells.Item(43, 10).Value = 19999  # J43
$ws.Cells.Item(43, 11).Value = 1990  # K43
$ws.Cells.Item(43, 12).Value = 19999  # L43
$ws.Cells.Item(43, 13).Value = -1839  # M43
$ws.Cells.Item(43, 14).Value = -20301  # N43
$ws.Cells.Item(76, 8).Value = 27057.666  # H76
$ws.Cells.Item(76, 10).Value = 27057.666  # J76
$ws.Cells.Item(76, 12).Value = 27057.666  # L76
$ws.Cells.Item(76, 14).Value = -27687.666  # N76
$ws.Cells.Item(79, 8).Value = 27057.666  # H79
$ws.Cells.Item(79, 10).Value = 27057.666  # J79
$ws.Cells.Item(79, 12).Value = 27057.666  # L79
$ws.Cells.Item(79, 14).Value = -29241.666  # N79

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 404.92856  # H55
$ws.Cells.Item(55, 9).Value = 245.90909  # I55
$ws.Cells.Item(55, 10).Value = 507.82352  # J55
$ws.Cells.Item(55, 11).Value = 245.90909  # K55
$ws.Cells.Item(55, 12).Value = 507.82352  # L55
$ws.Cells.Item(55, 13).Value = -72.90908999999999  # M55
$ws.Cells.Item(55, 14).Value = -853.8235199999999  # N55
$ws.Cells.Item(122, 8).Value = 1909.5333  # H122
$ws.Cells.Item(122, 9).Value = 1844.8334  # I122
$ws.Cells.Item(122, 10).Value = 2168.3333  # J122
$ws.Cells.Item(122, 11).Value = 5534.5002  # K122
$ws.Cells.Item(122, 12).Value = 6504.999899999999  # L122
$ws.Cells.Item(122, 13).Value = -3084.5002  # M122
$ws.Cells.Item(122, 14).Value = -11404.9999  # N122
$ws.Cells.Item(136, 8).Value = 3298  # H136
$ws.Cells.Item(136, 9).Value = 5700  # I136
$ws.Cells.Item(136, 11).Value = 17100  # K136
$ws.Cells.Item(136, 13).Value = -14550  # M136
$ws.Cells.Item(140, 8).Value = 57182.43  # H140
$ws.Cells.Item(140, 10).Value = 57182.43  # J140
$ws.Cells.Item(140, 12).Value = 57182.43  # L140
$ws.Cells.Item(140, 14).Value = -67542.42999999999  # N140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 21212.715  # H4
$ws.Cells.Item(4, 9).Value = 0  # I4
$ws.Cells.Item(4, 10).Value = 21212.715  # J4
$ws.Cells.Item(4, 11).Value = 0  # K4
$ws.Cells.Item(4, 12).Value = 21212.715  # L4
$ws.Cells.Item(4, 13).ClearContents()  # M4
$ws.Cells.Item(4, 14).Value = -21438.715  # N4
